$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2 down to whole numbers.
$ws.Range("Q2").Value = 612420
$ws.Range("R2").Value = 7034836

# Remove the empty/placeholder time cells (Starttid/Sluttid = "00:00").
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
